$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1473
$ws1.Range("F4").Value = 1130
$ws1.Range("F6").Value = 229
$ws1.Range("F8").Value = 688
$ws1.Range("F11").Value = 100
$ws1.Range("F12").Value = 225
$ws1.Range("F14").Value = 3242
$ws1.Range("F17").Value = 444
$ws1.Range("F19").Value = 513
$ws1.Range("F20").Value = 288
$ws1.Range("F24").Value = 680
$ws1.Range("F25").Value = 61
$ws1.Range("F26").Value = 257
$ws1.Range("F29").Value = 1608
$ws1.Range("F30").Value = 345

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 242
$ws2.Range("F7").Value = 238
$ws2.Range("F9").Value = 68

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 89

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1473
$ws4.Range("F5").Value = 1130
$ws4.Range("F8").Value = 90
$ws4.Range("F10").Value = 229
$ws4.Range("F12").Value = 688
$ws4.Range("F16").Value = 100
$ws4.Range("F17").Value = 225
$ws4.Range("F19").Value = 3242
$ws4.Range("F22").Value = 242
$ws4.Range("F23").Value = 444
$ws4.Range("F25").Value = 513
$ws4.Range("F26").Value = 288
$ws4.Range("F31").Value = 238
$ws4.Range("F33").Value = 68
$ws4.Range("F34").Value = 680
$ws4.Range("F38").Value = 61
$ws4.Range("F39").Value = 257
$ws4.Range("F42").Value = 1608
$ws4.Range("F43").Value = 345
